# Lab01 Review Report - fill in Coding Phase Defects sheet + small follow-on
# edits on the other two review sheets, matching the "Finished report for
# code analysis" commit.

$wb = $excel.ActiveWorkbook

$wsReq  = $wb.Worksheets.Item(1)   # Requirements Phase Defects
$wsArch = $wb.Worksheets.Item(2)   # Architect. Design Phase Defects
$wsCode = $wb.Worksheets.Item(3)   # Coding Phase Defects

$xlCenter = -4108

# ---------------------------------------------------------------------------
# 1) Requirements Phase Defects: effort-to-review cell becomes descriptive
#    text instead of a bare numeric hour count.
# ---------------------------------------------------------------------------
$wsReq.Range("E27").Value = "1 ora"

# ---------------------------------------------------------------------------
# 2) Architect. Design Phase Defects: effort-to-review cell becomes
#    descriptive text instead of the time-formatted fraction.
# ---------------------------------------------------------------------------
$wsArch.Range("E28").Value = "1 ora si 30 minute"

# ---------------------------------------------------------------------------
# 3) Coding Phase Defects: fill in the whole checklist table (rows 10-23)
#    with the checked item code, the doc page/line, and the
#    comments/improvements text, then center + middle-align that block
#    (rows 10-30) the way the reviewer formatted the finished report.
# ---------------------------------------------------------------------------

$codeRows = @(
    @(10, "C01", "Corigent line 117", "if ( media >= 4.5) trebuie transformat in if (media < 4.5)"),
    @(11, "C01", "Corigent line 110", "if(nrNote >= 0) nu este necesar"),
    @(12, "C01", "Corigent line 101", "if (clasa.size() >= 0) nu este necesar"),
    @(13, "C02", "ClasaRepositoryMock line 56 ", "if(clasa.size() >= 0) trebuie transformat in if (clasa.size() > 0) pentru a executa branch-ul pentru situatia in care clasa este goala"),
    @(14, "C03", $null, "Nu"),
    @(15, "C04", "main, linie 34", "Da, daca nu exista fisierul pentru studenti/ note"),
    @(16, "C05", $null, "Nu"),
    @(17, "C06", $null, "Nu"),
    @(18, "C07", $null, "Nu"),
    @(19, "C08", $null, "Nu"),
    @(20, "C09", $null, "Nu"),
    @(21, "C10", $null, "Nu"),
    @(22, "C11", $null, "Nu"),
    @(23, "C12", $null, "Nu")
)

foreach ($row in $codeRows) {
    $r = $row[0]
    $wsCode.Cells.Item($r, 3).Value = $row[1]
    if ($null -ne $row[2]) {
        $wsCode.Cells.Item($r, 4).Value = $row[2]
    }
    $wsCode.Cells.Item($r, 5).Value = $row[3]
}

# Row heights for the entries whose comment text wraps onto multiple lines.
$wsCode.Rows.Item(10).RowHeight = 30
$wsCode.Rows.Item(13).RowHeight = 45
$wsCode.Rows.Item(15).RowHeight = 30

# Center + middle-align the whole answered block (rows 10-30), matching the
# reviewer's final formatting pass.
$block = $wsCode.Range("C10:E30")
$block.HorizontalAlignment = $xlCenter
$block.VerticalAlignment = $xlCenter

# Trailing spacer row gains the same centered styling, now stretching one
# column further (through F) because of the widened used range below, and
# the final effort-to-review cell gets its descriptive text.
$spacer = $wsCode.Range("C31:F31")
$spacer.HorizontalAlignment = $xlCenter
$spacer.VerticalAlignment = $xlCenter

$wsCode.Range("E32").Value = "2 ore"

# ---------------------------------------------------------------------------
# 4) Selections left by the reviewer when the workbook was last saved.
#    Apply the non-active sheets first so the final Activate() below leaves
#    the "Architect. Design Phase Defects" tab selected, as in the source.
# ---------------------------------------------------------------------------
$wsReq.Range("K24").Select()
$wsCode.Range("E32").Select()

$wsArch.Activate()
$wsArch.Range("E29").Select()
